$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "customer service" title text in B5 to the new page title.
$ws.Range("B5").Value = "Help & Contact Us - Amazon Customer Service"

# Add a hyperlink on B2 (the "amazon home page url" cell) pointing at the
# Amazon home page, mirroring the existing hyperlink pattern on B6.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.amazon.com/")

# Move the active selection to B5.
$ws.Range("B5").Select() | Out-Null
